$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update coin prices (Price column, D) to the latest scraped values.
# Values are written as text (apostrophe-prefixed) to preserve the original
# inline-string cell type used by this sheet, then the quote-prefix style
# flag is cleared so no stray formatting is introduced.
$ws.Range("D2").Value = "'267.72"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'21.36"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'6.246"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'0.06204"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'3.569"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'6.541"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'1.386"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.8251"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'0.1633"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'0.08271"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'0.03558"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'0.03183"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'0.09198"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'3.764"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'0.001629"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'0.04682"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").Value = "'0.006205"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'0.001067"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "'0.0001500"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").Value = "'2.289"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = "'0.01367"
$ws.Range("D24").Style = "Normal"
$ws.Range("D28").Value = "'0.0002713"
$ws.Range("D28").Style = "Normal"
$ws.Range("D40").Value = "'0.04712"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = "'0.006975"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Value = "'0.004200"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Value = "'0.1119"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "'0.01149"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.00006201"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "'0.0009900"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "'0.8025"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Value = "'0.002340"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Value = "'0.00001900"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Value = "'0.01240"
$ws.Range("D51").Style = "Normal"

# Row 48/49 "Worstin24h" suffix moved from BOLO (49) to CoinbaseStockToken (48)
# following the updated 24h ranking.
$ws.Range("E48").Value = "47CoinbaseStockTokenCOINWorstin24h"
$ws.Range("E49").Value = "48BOLOBOLO"
